$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width: 25 -> 28 characters (raw OOXML width units) ---
# Excel's ColumnWidth COM property uses a slightly different unit than the
# raw stored width; 27.14 round-trips to a stored width of exactly 28.
$ws.Columns.Item(1).ColumnWidth = 27.14

# --- New "Abfragen Umlandregionen" block (rows 17-21), column D examples ---
$ws.Range("D17").Value = "alle Umzüge aus Oberzentren MV in ABL"
$ws.Range("D17").WrapText = $true

$ws.Range("D18").Value = "alle Umzüge aus Umlandregionen MV in ABL"
$ws.Range("D18").WrapText = $true

$ws.Range("D19").Value = "alle Umzüge aus Umlandregionen von Rostock nach Rostock"
$ws.Range("D19").WrapText = $true

$ws.Range("D20").Value = "Definition der Umlandregionen aus Pendlerrate zu Oberzentren(Zone1>= 40%; zone2>=20%) - Pendlerdaten 2006 benötigt"
$ws.Range("D20").WrapText = $true

$ws.Range("D21").Value = "?"
$ws.Range("D21").WrapText = $true

# --- "Abfragen Kreisdaten" block (rows 23-32), column D examples ---
$ws.Range("D23").Value = "alle Fortzüge Kreis Potsdam-Mittelmark in Kreise ABL 2004"
$ws.Range("D24").Value = "alle Zuzüge Kreis Potsdam-Mittelmark aus Kreisen ABL 2004"
$ws.Range("D25").Value = "alle Zuzüge Kreis Potsdam-Mittelmark aus Kreisen ABL 2004 G:w, FS:1"
$ws.Range("D26").Value = "Bevölkerungssaldo Sachsen-Anhalt auf Kreisebene 2000 bis 2006"
$ws.Range("D27").Value = "Wanderungsssaldo Sachsen-Anhalt auf Kreisebene 2000 bis 2006"
$ws.Range("D28").Value = "Wanderungsssaldo Sachsen-Anhalt auf Kreisebene 2000 bis 2006/EW"
$ws.Range("D29").Value = "?Definition Umlandkreise?"
$ws.Range("D30").Value = "alle Kreise BB mit Wanderungssaldo/EW von >2%"

$ws.Range("D31").Value = "alle Zuzüge auf Kreisebene von ABL zu NBL 2006"
$ws.Range("D31").WrapText = $true

$ws.Range("D32").Value = "?Definition Umlandkreise?"

# --- "Abfragen Bundesländer" block (rows 34-39), column D examples ---
$ws.Range("D34").Value = "alle Fortzüge Thüringen in ABL 2003"
$ws.Range("D35").Value = "alle Zuzüge Thüringen aus ABL 2003"
$ws.Range("D36").Value = "alle Fortzüge Thüringen in ABL 2003 G:w, FS:1"
$ws.Range("D37").Value = "Bevölkerungssaldos BL"
$ws.Range("D38").Value = "Wanderungssaldos NBL"
$ws.Range("D39").Value = "Wanderungssaldos NBL/EW"

# --- New "allgemeine, einfache Abfragen" block (rows 41-48) ---
$ws.Range("A41").Value = "allgemeine, einfache  Abfragen"
$ws.Range("D41").Value = "Bevölkerungsstand 2000 bis 2009 Deutschland"
$ws.Range("D42").Value = "Bevölkerungsstand 2000 bis 2009 Sachsen"
$ws.Range("D43").Value = "Bevölkerungsstand 2000 bis 2009 Gemeinde (14729)"
$ws.Range("D44").Value = "Bevölkerungsstand 2000 bis 2009 Gemeinde Wurzen, Stadt (14729410)"
$ws.Range("D45").Value = "Bevölkerungsstand D 2000 nach BL"
$ws.Range("D46").Value = "Bevölkerungsstand D 2000 nach Kreisen"
$ws.Range("D47").Value = "Bevölkerungsstand D 2000 nach Gemeinden"
$ws.Range("D48").Value = "Umzüge D 2000 nach NBL"

# --- Update the active selection to match the author's final cursor position ---
$ws.Range("D49").Select()
